{"js": "// Map of old text -> new text, taken from the diff.\nconst replacements = [\n  [\"2024-03-20 Wednesday\", \"2024-03-21 Thursday\"],\n  [\"18\u00d718=\", \"23\u00d746=\"],\n  [\"95\u00d788=\", \"44\u00d785=\"],\n  [\"14\u00d783=\", \"49\u00d715=\"],\n  [\"65\u00d796=\", \"48\u00d748=\"],\n  [\"87\u00d760=\", \"60\u00d762=\"],\n  [\"71\u00d756=\", \"41\u00d724=\"],\n  [\"15\u00d757=\", \"34\u00d736=\"],\n  [\"76\u00d734=\", \"43\u00d782=\"],\n  [\"72\u00d763=\", \"78\u00d791=\"],\n  [\"60\u00d711=\", \"23\u00d752=\"],\n  [\"68\u00d774=\", \"47\u00d761=\"],\n  [\"63\u00d719=\", \"19\u00d783=\"],\n  [\"66\u00d729=\", \"39\u00d787=\"],\n  [\"49\u00d774=\", \"58\u00d775=\"],\n  [\"95\u00d756=\", \"36\u00d713=\"],\n  [\"42\u00d749=\", \"28\u00d785=\"],\n  [\"71\u00d722=\", \"19\u00d742=\"],\n  [\"91\u00d796=\", \"37\u00d766=\"],\n  [\"98\u00d778=\", \"11\u00d787=\"],\n  [\"26\u00d759=\", \"80\u00d786=\"],\n  [\"78\u00d743=\", \"48\u00d713=\"],\n  [\"58\u00d770=\", \"68\u00d714=\"],\n  [\"51\u00d797=\", \"56\u00d720=\"],\n  [\"37\u00d734=\", \"68\u00d760=\"],\n  [\"37\u00d755=\", \"47\u00d735=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Map of old text -> new text, taken from the diff.\n$replacements = @(\n    @(\"2024-03-20 Wednesday\", \"2024-03-21 Thursday\"),\n    @(\"18\u00d718=\", \"23\u00d746=\"),\n    @(\"95\u00d788=\", \"44\u00d785=\"),\n    @(\"14\u00d783=\", \"49\u00d715=\"),\n    @(\"65\u00d796=\", \"48\u00d748=\"),\n    @(\"87\u00d760=\", \"60\u00d762=\"),\n    @(\"71\u00d756=\", \"41\u00d724=\"),\n    @(\"15\u00d757=\", \"34\u00d736=\"),\n    @(\"76\u00d734=\", \"43\u00d782=\"),\n    @(\"72\u00d763=\", \"78\u00d791=\"),\n    @(\"60\u00d711=\", \"23\u00d752=\"),\n    @(\"68\u00d774=\", \"47\u00d761=\"),\n    @(\"63\u00d719=\", \"19\u00d783=\"),\n    @(\"66\u00d729=\", \"39\u00d787=\"),\n    @(\"49\u00d774=\", \"58\u00d775=\"),\n    @(\"95\u00d756=\", \"36\u00d713=\"),\n    @(\"42\u00d749=\", \"28\u00d785=\"),\n    @(\"71\u00d722=\", \"19\u00d742=\"),\n    @(\"91\u00d796=\", \"37\u00d766=\"),\n    @(\"98\u00d778=\", \"11\u00d787=\"),\n    @(\"26\u00d759=\", \"80\u00d786=\"),\n    @(\"78\u00d743=\", \"48\u00d713=\"),\n    @(\"58\u00d770=\", \"68\u00d714=\"),\n    @(\"51\u00d797=\", \"56\u00d720=\"),\n    @(\"37\u00d734=\", \"68\u00d760=\"),\n    @(\"37\u00d755=\", \"47\u00d735=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
